# "3/4 of robustness check"
# Re-labels the company axis (ไทยเดนมาร์ค -> ดัชมิลล์ swapped in, "หนองโพ" dropped,
# remaining companies re-ordered) and writes the freshly recomputed Jaccard
# similarity matrix (self-similarity now correctly reads 1 on the diagonal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row/column label order -------------------------------------------------
$labels = @("เมจิ", "ดัชมิลล์", "โฟร์โมสต์", "แดรี่โฮม", "โชคชัย", "เอ็มมิลค์", "ไทยเดนมาร์ค")

for ($i = 0; $i -lt $labels.Length; $i++) {
    # Header row: B1:H1
    $ws.Cells.Item(1, $i + 2).Value = $labels[$i]
    # Label column: A2:A8
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# --- New similarity matrix values (B2:H8) ----------------------------------------
# (PowerShell-style numeric literals here avoid scientific notation - the engine's
# parser only accepts plain decimals.)
$matrix = @(
    @(1.0, 0.1108032196764444, 0.09446707970288731, 0.01956254241693767, 0.02964585807536208, 0.04017087007589609, 0.05003973575996431),
    @(0.1108032196764444, 1.0, 0.124671961322111, 0.02189330424219897, 0.04173466467630537, 0.0004553734061930783, 0.1114107448998216),
    @(0.09446707970288731, 0.124671961322111, 1.0, 0, 0.03547041304011841, 0.0002579979360165118, 0.117420824480639),
    @(0.01956254241693767, 0.02189330424219897, 0, 1.0, 0.01346409356032294, 0.003787878787878788, 0.0426602026144021),
    @(0.02964585807536208, 0.04173466467630537, 0.03547041304011841, 0.01346409356032294, 1.0, 0, 0.02487639490014739),
    @(0.04017087007589609, 0.0004553734061930783, 0.0002579979360165118, 0.003787878787878788, 0, 1.0, 0.00004782629489693433),
    @(0.05003973575996431, 0.1114107448998216, 0.117420824480639, 0.0426602026144021, 0.02487639490014739, 0.00004782629489693433, 1.0)
)

for ($r = 0; $r -lt 7; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $matrix[$r][$c]
    }
}
